$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Wipe existing body content (keep header row labels A1/B1) ---
$ws.Range("C1:L15").Clear()
$ws.Range("A2:L15").Clear()

# --- Row 1 (header row totals) ---
$ws.Range("C1").Value = 395
$ws.Range("D1").Value = 218

# --- Data rows (empresa - 2 import) ---
$ws.Range("B2").Value = 123
$ws.Range("C2").Value = "53:17"

$ws.Range("B3").Value = 130
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = 325

$ws.Range("B4").Value = 83
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = 299

$ws.Range("B5").Value = 134
$ws.Range("C5").Value = "07:07"
$ws.Range("D5").Value = 208

$ws.Range("B6").Value = 105
$ws.Range("C6").Value = "97:52"

$ws.Range("B7").Value = 118
$ws.Range("C7").Value = "60:20"

$ws.Range("B8").Value = 78
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = 338

$ws.Range("B9").Value = 135
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = 325

$ws.Range("B10").Value = 117
$ws.Range("C10").Value = ""

$ws.Range("B11").Value = 126
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = 325

$ws.Range("B12").Value = ""
$ws.Range("C12").Value = "218:38"

# trailing near-blank rows
$ws.Range("B13").Style = "Normal"
$ws.Range("B13").Font.Bold = $false
$ws.Range("B14").Style = "Normal"
$ws.Range("B14").Font.Bold = $false

# --- Formatting: box border + right-align for the data block ---
$dataRange = $ws.Range("B2:D12")
$dataRange.Borders.LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$dataRange.Borders.Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin
$dataRange.Borders.ColorIndex = [Microsoft.Office.Interop.Excel.XlColorIndex]::xlColorIndexAutomatic
$dataRange.HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignRight

# --- Helper cell (bold font, [h]:mm:ss format), placed so that after the
#     column delete below it lands on I6 ---
$ws.Range("L6").NumberFormat = "[h]:mm:ss"
$ws.Range("L6").Font.Bold = $true

# --- Remove now-unused columns E:G (their bestFit width no longer needed) ---
$ws.Columns("E:G").Delete()

# --- Update selection to match author's final cursor position ---
$ws.Range("D20").Select()
